$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in T3:W3 from 0.99 to 1
$ws.Range("T3:W3").Value = 1

# Update the selected range on the sheet view to match the actual used dimension
$ws.Range("A1:X7").Select()
